# Generate Report for Handback
# Update timestamps / priority values produced by a later report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# bc7d2a97-... row and eda47f45-... row both updated to the new timestamp.
$wsOverview.Range("G4").Value = "2016-08-30 02:17:08"
$wsOverview.Range("G5").Value = "2016-08-30 02:17:08"

# zh-cn sheet: Priority (E), Correspond Handoff Datetime (H), and
# Correspond Handback DateTime (K) for rows 4 and 5.
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-30 02:16:58"
$wsZhCn.Range("H5").Value = "2016-08-30 02:16:58"
$wsZhCn.Range("K4").Value = "2016-08-30 02:17:51"
$wsZhCn.Range("K5").Value = "2016-08-30 02:17:51"

# de-de sheet: Priority (E), Correspond Handoff Datetime (H), and
# Correspond Handback DateTime (K) for rows 4 and 5.
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H4").Value = "2016-08-30 02:17:08"
$wsDeDe.Range("H5").Value = "2016-08-30 02:17:08"
$wsDeDe.Range("K4").Value = "2016-08-30 02:17:58"
$wsDeDe.Range("K5").Value = "2016-08-30 02:17:58"
